$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the new column (CA / col 79) width to the existing BZ (col 78) column width
$refColWidth = $ws.Columns.Item(78).ColumnWidth()
$ws.Columns.Item(79).ColumnWidth = $refColWidth

# Header cell CA1: new date label "2024/11/26", stored as text (matching the other date headers)
$ws.Range("CA1").NumberFormat = "@"
$ws.Range("CA1").Value = "2024/11/26"
$ws.Range("BZ1").Copy()
$ws.Range("CA1").PasteSpecial(-4122)

# Data cells CA2:CA53 - value + style copied from a same-styled existing cell
# (style 1 = default, style 2 = yellow highlight for <125, style 3 = light blue for 125-139.9)
$ws.Range("A2").Copy()
$ws.Range("CA2").PasteSpecial(-4122)
$ws.Range("CA2").Value = 142.3
$ws.Range("A2").Copy()
$ws.Range("CA3").PasteSpecial(-4122)
$ws.Range("CA3").Value = 149.5
$ws.Range("A2").Copy()
$ws.Range("CA4").PasteSpecial(-4122)
$ws.Range("CA4").Value = 231.3
$ws.Range("D2").Copy()
$ws.Range("CA5").PasteSpecial(-4122)
$ws.Range("CA5").Value = 107.4
$ws.Range("A2").Copy()
$ws.Range("CA6").PasteSpecial(-4122)
$ws.Range("CA6").Value = 172.7
$ws.Range("A2").Copy()
$ws.Range("CA7").PasteSpecial(-4122)
$ws.Range("CA7").Value = 166.2
$ws.Range("A2").Copy()
$ws.Range("CA8").PasteSpecial(-4122)
$ws.Range("CA8").Value = 220.2
$ws.Range("D2").Copy()
$ws.Range("CA9").PasteSpecial(-4122)
$ws.Range("CA9").Value = 124.5
$ws.Range("N2").Copy()
$ws.Range("CA10").PasteSpecial(-4122)
$ws.Range("CA10").Value = 139.2
$ws.Range("A2").Copy()
$ws.Range("CA11").PasteSpecial(-4122)
$ws.Range("CA11").Value = 179.1
$ws.Range("N2").Copy()
$ws.Range("CA12").PasteSpecial(-4122)
$ws.Range("CA12").Value = 134
$ws.Range("A2").Copy()
$ws.Range("CA13").PasteSpecial(-4122)
$ws.Range("CA13").Value = 188.8
$ws.Range("N2").Copy()
$ws.Range("CA14").PasteSpecial(-4122)
$ws.Range("CA14").Value = 134.8
$ws.Range("A2").Copy()
$ws.Range("CA15").PasteSpecial(-4122)
$ws.Range("CA15").Value = 152.9
$ws.Range("N2").Copy()
$ws.Range("CA16").PasteSpecial(-4122)
$ws.Range("CA16").Value = 127.9
$ws.Range("A2").Copy()
$ws.Range("CA17").PasteSpecial(-4122)
$ws.Range("CA17").Value = 176.3
$ws.Range("N2").Copy()
$ws.Range("CA18").PasteSpecial(-4122)
$ws.Range("CA18").Value = 133.4
$ws.Range("A2").Copy()
$ws.Range("CA19").PasteSpecial(-4122)
$ws.Range("CA19").Value = 170.4
$ws.Range("N2").Copy()
$ws.Range("CA20").PasteSpecial(-4122)
$ws.Range("CA20").Value = 134
$ws.Range("A2").Copy()
$ws.Range("CA21").PasteSpecial(-4122)
$ws.Range("CA21").Value = 263.6
$ws.Range("A2").Copy()
$ws.Range("CA22").PasteSpecial(-4122)
$ws.Range("CA22").Value = 234.3
$ws.Range("A2").Copy()
$ws.Range("CA23").PasteSpecial(-4122)
$ws.Range("CA23").Value = 151.3
$ws.Range("D2").Copy()
$ws.Range("CA24").PasteSpecial(-4122)
$ws.Range("CA24").Value = 121.8
$ws.Range("A2").Copy()
$ws.Range("CA25").PasteSpecial(-4122)
$ws.Range("CA25").Value = 207.5
$ws.Range("A2").Copy()
$ws.Range("CA26").PasteSpecial(-4122)
$ws.Range("CA26").Value = 142
$ws.Range("N2").Copy()
$ws.Range("CA27").PasteSpecial(-4122)
$ws.Range("CA27").Value = 128.3
$ws.Range("A2").Copy()
$ws.Range("CA28").PasteSpecial(-4122)
$ws.Range("CA28").Value = 203.7
$ws.Range("A2").Copy()
$ws.Range("CA29").PasteSpecial(-4122)
$ws.Range("CA29").Value = 144.4
$ws.Range("A2").Copy()
$ws.Range("CA30").PasteSpecial(-4122)
$ws.Range("CA30").Value = 149.9
$ws.Range("A2").Copy()
$ws.Range("CA31").PasteSpecial(-4122)
$ws.Range("CA31").Value = 301.3
$ws.Range("A2").Copy()
$ws.Range("CA32").PasteSpecial(-4122)
$ws.Range("CA32").Value = 146.4
$ws.Range("A2").Copy()
$ws.Range("CA33").PasteSpecial(-4122)
$ws.Range("CA33").Value = 208.4
$ws.Range("A2").Copy()
$ws.Range("CA34").PasteSpecial(-4122)
$ws.Range("CA34").Value = 144.2
$ws.Range("N2").Copy()
$ws.Range("CA35").PasteSpecial(-4122)
$ws.Range("CA35").Value = 135.3
$ws.Range("A2").Copy()
$ws.Range("CA36").PasteSpecial(-4122)
$ws.Range("CA36").Value = 166.7
$ws.Range("A2").Copy()
$ws.Range("CA37").PasteSpecial(-4122)
$ws.Range("CA37").Value = 188.9
$ws.Range("A2").Copy()
$ws.Range("CA38").PasteSpecial(-4122)
$ws.Range("CA38").Value = 171.5
$ws.Range("A2").Copy()
$ws.Range("CA39").PasteSpecial(-4122)
$ws.Range("CA39").Value = 174.4
$ws.Range("A2").Copy()
$ws.Range("CA40").PasteSpecial(-4122)
$ws.Range("CA40").Value = 144.6
$ws.Range("A2").Copy()
$ws.Range("CA41").PasteSpecial(-4122)
$ws.Range("CA41").Value = 242.7
$ws.Range("A2").Copy()
$ws.Range("CA42").PasteSpecial(-4122)
$ws.Range("CA42").Value = 171.2
$ws.Range("A2").Copy()
$ws.Range("CA43").PasteSpecial(-4122)
$ws.Range("CA43").Value = 186.6
$ws.Range("N2").Copy()
$ws.Range("CA44").PasteSpecial(-4122)
$ws.Range("CA44").Value = 132.5
$ws.Range("A2").Copy()
$ws.Range("CA45").PasteSpecial(-4122)
$ws.Range("CA45").Value = 252.6
$ws.Range("A2").Copy()
$ws.Range("CA46").PasteSpecial(-4122)
$ws.Range("CA46").Value = 149.3
$ws.Range("A2").Copy()
$ws.Range("CA47").PasteSpecial(-4122)
$ws.Range("CA47").Value = 184.4
$ws.Range("A2").Copy()
$ws.Range("CA48").PasteSpecial(-4122)
$ws.Range("CA48").Value = 168.7
$ws.Range("A2").Copy()
$ws.Range("CA49").PasteSpecial(-4122)
$ws.Range("CA49").Value = 353.7
$ws.Range("N2").Copy()
$ws.Range("CA50").PasteSpecial(-4122)
$ws.Range("CA50").Value = 131.9
$ws.Range("A2").Copy()
$ws.Range("CA51").PasteSpecial(-4122)
$ws.Range("CA51").Value = 199.2
$ws.Range("A2").Copy()
$ws.Range("CA52").PasteSpecial(-4122)
$ws.Range("CA52").Value = 165.2
$ws.Range("A2").Copy()
$ws.Range("CA53").PasteSpecial(-4122)
$ws.Range("CA53").Value = 167.7
